$d = $word.ActiveDocument

$pairs = @(
    @("521÷8=", "787÷4="),
    @("535÷3=", "339÷6="),
    @("836÷3=", "745÷4="),
    @("186÷5=", "119÷3="),
    @("346÷4=", "784÷2="),
    @("802÷9=", "391÷2="),
    @("216÷9=", "179÷9="),
    @("776÷6=", "323÷2="),
    @("410÷4=", "256÷6="),
    @("474÷9=", "267÷4="),
    @("138÷2=", "260÷7="),
    @("767÷8=", "540÷2="),
    @("621÷2=", "678÷4="),
    @("310÷2=", "301÷9="),
    @("746÷9=", "369÷3="),
    @("641÷5=", "450÷8="),
    @("748÷9=", "308÷8="),
    @("940÷9=", "314÷2="),
    @("394÷7=", "120÷6="),
    @("688÷5=", "558÷5="),
    @("847÷9=", "592÷9="),
    @("816÷3=", "448÷4="),
    @("843÷6=", "577÷7="),
    @("710÷7=", "673÷2="),
    @("529÷7=", "889÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
